$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove row 577 ("「テレビを観るのをやめて下さい」") and shift all
# subsequent rows up by one, matching the target diff.
$ws.Rows.Item(577).Delete()
